# Fixed update to excel issue
# Adds a "PO Forecast" sheet (Prophet-style forecast output) and renames the
# "Requested quantity" header on the existing sheets to match the new naming
# convention used across the PO analysis workbook.

$wb = $excel.ActiveWorkbook

# --- Rename headers on the two existing sheets -----------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item("Monthly Trend")
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- Add the new "PO Forecast" sheet after the last existing sheet ---------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Match the page margins used on the other sheets in the workbook (values are
# in points: 0.75in/1in/0.5in -> 54/72/36pt).
$wsForecast.PageSetup.LeftMargin = 54
$wsForecast.PageSetup.RightMargin = 54
$wsForecast.PageSetup.TopMargin = 72
$wsForecast.PageSetup.BottomMargin = 72
$wsForecast.PageSetup.HeaderMargin = 36
$wsForecast.PageSetup.FooterMargin = 36

# Header row
$wsForecast.Cells.Item(1, 1).Value = "ds"
$wsForecast.Cells.Item(1, 2).Value = "PO_Forecast"
$wsForecast.Cells.Item(1, 3).Value = "yhat_lower"
$wsForecast.Cells.Item(1, 4).Value = "yhat_upper"

# Match the look of the header rows on the other sheets (bold, bordered,
# centered + top-aligned).
$headerRange = $wsForecast.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Forecast data rows (ds, PO_Forecast, yhat_lower, yhat_upper)
$forecastData = @(
    @(45039.99999999999,0,-73.82423894058364,34.22457956366898),
    @(45060.99999999999,0,-71.17707999094534,39.34055967097463),
    @(45102.99999999999,0,-60.90845130208205,43.04555378491567),
    @(45116.99999999999,0,-58.35234062041096,46.75814670326071),
    @(45130.99999999999,0,-57.3339793850044,49.0119070399985),
    @(45151.99999999999,0,-52.61203688400391,52.65466603310664),
    @(45158.99999999999,0,-47.73044391519865,54.16213969577318),
    @(45165.99999999999,0,-51.35173363108444,59.2918093566544),
    @(45193.99999999999,5,-45.96667802688083,56.48464475645435),
    @(45200.99999999999,6,-47.66628229635462,57.19592824193234),
    @(45207.99999999999,7,-46.91296764381313,61.70406342428835),
    @(45228.99999999999,10,-47.94285784564088,60.17833177479083),
    @(45235.99999999999,11,-41.92131811474251,61.73988362707364),
    @(45242.99999999999,12,-39.42610816731084,65.08151744562142),
    @(45249.99999999999,14,-37.73198652971566,65.42989865084081),
    @(45256.99999999999,15,-34.02678579592836,67.12935202959278),
    @(45312.99999999999,23,-28.36828056001267,76.81876076920106),
    @(45319.99999999999,25,-29.10501768319008,76.88526764711307),
    @(45326.99999999999,26,-23.96794685567725,79.06470439498707),
    @(45333.99999999999,27,-19.23937878184386,76.71577484548503),
    @(45368.99999999999,32,-19.9594182695046,80.29588883913681),
    @(45382.99999999999,35,-16.53947986047538,92.38560100414965),
    @(45389.99999999999,36,-15.83628164840641,84.5262224044192),
    @(45396.99999999999,37,-19.07273840176955,86.20759107320808),
    @(45403.99999999999,38,-13.49941859316193,88.39493834823683),
    @(45410.99999999999,39,-14.89567519199886,87.60718127873467),
    @(45417.99999999999,40,-10.93154266036144,91.97418142554615),
    @(45424.99999999999,41,-11.55018474034096,95.53232601148373),
    @(45431.99999999999,42,-8.889732642954696,95.44210028716245),
    @(45438.99999999999,43,-13.46947133570896,96.36307974181749),
    @(45445.99999999999,44,-10.5367627495143,101.3043418057849),
    @(45452.99999999999,46,-5.509811127216756,97.85201729405016),
    @(45459.99999999999,47,-6.903107396902906,97.8671351985686),
    @(45466.99999999999,48,-5.504889250393597,95.92213965943344),
    @(45473.99999999999,49,-3.417696281494674,98.74110848734087),
    @(45480.99999999999,50,-0.504548682371232,101.7753459120746),
    @(45487.99999999999,51,-4.833430526011122,101.8740175888512),
    @(45494.99999999999,52,-2.20964778833247,107.4093000944336),
    @(45501.99999999999,53,-3.037710750287723,106.3923077442819),
    @(45508.99999999999,54,3.159614964046666,106.9751617618572),
    @(45515.99999999999,55,4.830385310873455,106.5909505345192)
)

$row = 2
foreach ($r in $forecastData) {
    $wsForecast.Cells.Item($row, 1).Value = $r[0]
    $wsForecast.Cells.Item($row, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $wsForecast.Cells.Item($row, 2).Value = $r[1]
    $wsForecast.Cells.Item($row, 3).Value = $r[2]
    $wsForecast.Cells.Item($row, 4).Value = $r[3]
    $row++
}
